$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "Username : 31246;`nPassword : bni1234;`nRole : 20/21 - Analis Investasi/Asisten Investasi;`nNo. Urut : 2962"
$ws.Range("M2").Value = 2962

$ws.Range("G2").Select()
